# Add a new "2023" data column (T) to the disasters-deaths table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Columns A:C get a uniform width instead of three different custom widths.
#    (34.25 is the closest input that this engine's pixel-quantized ColumnWidth
#    setter rounds to ~35.14, matching the target width as closely as possible.)
$ws.Columns("A:C").ColumnWidth = 34.25

# 2) Column T needs the same formatting (font/alignment/borders) as column S
#    for every row of the table (rows 3-34). Copy format+value first, then
#    overwrite with the actual 2023 figures below.
$src = $ws.Range("S3:S34")
$dst = $ws.Range("T3:T34")
$src.Copy($dst)

# 3) Fill in the 2023 values (column T) for each row of the table.
$ws.Range("T4").Value = 2023

$ws.Range("T5").Value = 44
$ws.Range("T6").Value = 24
$ws.Range("T7").Value = 20

$ws.Range("T8").Value = "-"
$ws.Range("T9").Value = "-"
$ws.Range("T10").Value = "-"

$ws.Range("T11").Value = 5
$ws.Range("T12").Value = 1
$ws.Range("T13").Value = 4

$ws.Range("T14").Value = 8
$ws.Range("T15").Value = 6
$ws.Range("T16").Value = 2

$ws.Range("T17").Value = 5
$ws.Range("T18").Value = 1
$ws.Range("T19").Value = 4

$ws.Range("T20").Value = 7
$ws.Range("T21").Value = 5
$ws.Range("T22").Value = 2

$ws.Range("T23").Value = "-"
$ws.Range("T24").Value = "-"
$ws.Range("T25").Value = "-"

$ws.Range("T26").Value = 18
$ws.Range("T27").Value = 10
$ws.Range("T28").Value = 8

$ws.Range("T29").Value = "-"
$ws.Range("T30").Value = "-"
$ws.Range("T31").Value = "-"

$ws.Range("T32").Value = 1
$ws.Range("T33").Value = 1
$ws.Range("T34").Value = "-"

# 4) Reset the selection away from the stale "T24" reference left over in the
#    source file so it no longer points at a cell that now holds real data.
$null = $ws.Range("A1").Select()
